$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'track tights',
    'wintergear compression men',
    'basketball training tights',
    'basketball protector',
    'winter leggings for men',
    'kids tights with knee pads',
    'ropa de monta?a hombre',
    'men workout leggings nike',
    'boys white knee pads basketball',
    'coyote brown pants with knee pads',
    'black basketball knee pads nike',
    'basketball knee pads womens',
    'cold weather workout pants for men',
    'mens under armour long underwear pants',
    'calf compression leggings men',
    'tights for men nike',
    'swimming pants for men',
    'copper compression pants for men',
    'tesla thermal pants',
    'hex knee sleeve',
    'thermal nike',
    'running base layer men',
    'long spandex men',
    'cold gear compression pants men',
    'soccer winter gear',
    'kids soccor gear',
    'mens pants with knee pads',
    'mens thermal baselayer',
    'knee pad pants men',
    'adidas compression pants men',
    'adult clothing protector',
    'adult football knee pads',
    'adult football pants with pads',
    'athletic knee pads',
    'athletic tights men',
    'baseball knee pad',
    'baseball knee pads',
    'baseball pants adult small',
    'baskerball tights',
    'basketball compression knee pads',
    'basketball compression pants youth with knee pads',
    'basketball for youth',
    'basketball hip pads',
    'basketball leggings men',
    'basketball pants for women',
    'basketball pants men',
    'basketball tights boys youth',
    'basketball tights for men mcdavid',
    'basketball tights with pads for boys',
    'best basketball knee pads',
    'big knee pads',
    'big man knee pads',
    'bjj leggings',
    'black knee pads',
    'black knee pads for volleyball',
    'black leggings design',
    'black mesh leggings capri',
    'black youth baseball pants',
    'boys athletic tights basketball',
    'boys basketball knee pads mcdavid',
    'boys compression leggings',
    'boys compression leggings youth',
    'boys compression tights',
    'break away basketball pants',
    'capri leggings medium',
    'capri mens',
    'capri tights for men',
    'capris leggings',
    'cold gear for football',
    'cold weather panta',
    'compression for knee',
    'compression gear',
    'compression men pants',
    'compression pants big and tall men',
    'compression pants padded knees basketball',
    'compression shorts 3 4 length men',
    'compression tights for men',
    'cycling pants',
    'dark purple basketball knee pads',
    'elbow knee pads youth',
    'excersize equipment for men',
    'football 3 4 tights',
    'football knee pads',
    'football leg pads',
    'football pants',
    'football pants adult black',
    'football pants youth',
    'g form knee pads youth',
    'gel knee pads',
    'gel knee pads for men',
    'girl knee pads',
    'girls basketball knee pads',
    'girls tights with knee pads',
    'girls volleyball knee pads',
    'girls volleyball knee pads youth',
    'girls youth volleyball knee pads',
    'gym pants for men',
    'happy knees',
    'hex foam knee pads',
    'hex knee pads for basketball'
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
